# Generate Report for Handback
#
# This mirrors a "Generate handback report" run: the e93ba108 row (row 6)
# in the Overview sheet and in each locale sheet (zh-cn, de-de) moves from
# "Ready for handoff" to "Handed back: in sync with en-US", and the
# per-locale handback columns (Latest Target File / Latest Handback File /
# Latest Handback DateTime) get filled in for that row, same as the other
# already-handed-back rows above it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 6 (e93ba108-6a4a-498d-84aa-cffefbf1ccb9.md)
# zh-cn / de-de status columns flip from "Ready for handoff" to
# "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E6").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F6").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet: row 6 status + handback columns
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C6").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("I6").Value = "e93ba108-6a4a-498d-84aa-cffefbf1ccb9.md"
$wsZhCn.Range("I6").Font.Underline = $true
$wsZhCn.Range("I6").Font.Color = 15570276
$wsZhCn.Range("J6").Value = "e93ba108-6a4a-498d-84aa-cffefbf1ccb9.dbd3997b3aa8b1c92fb232b2fbd9be9e105c8226.zh-cn.xlf"
$wsZhCn.Range("K6").Value = "2016-09-07 06:02:12"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/03487e80efda313ca52c5331e6f5e6dc67100d41/e2e/e93ba108-6a4a-498d-84aa-cffefbf1ccb9.md", $null, $null, "e93ba108-6a4a-498d-84aa-cffefbf1ccb9.md") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet: row 6 status + handback columns
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C6").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("I6").Value = "e93ba108-6a4a-498d-84aa-cffefbf1ccb9.md"
$wsDeDe.Range("I6").Font.Underline = $true
$wsDeDe.Range("I6").Font.Color = 15570276
$wsDeDe.Range("J6").Value = "e93ba108-6a4a-498d-84aa-cffefbf1ccb9.dbd3997b3aa8b1c92fb232b2fbd9be9e105c8226.de-de.xlf"
$wsDeDe.Range("K6").Value = "2016-09-07 06:02:30"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/03487e80efda313ca52c5331e6f5e6dc67100d41/e2e/e93ba108-6a4a-498d-84aa-cffefbf1ccb9.md", $null, $null, "e93ba108-6a4a-498d-84aa-cffefbf1ccb9.md") | Out-Null

Write-Output "Handback report generated for e93ba108-6a4a-498d-84aa-cffefbf1ccb9.md"
